$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -ne $null) {
        $text = $val.ToString()
        if ($text.StartsWith("System, ")) {
            $rest = $text.Substring(8)
            $cell.Value = "$rest, System"
        }
    }
}
